$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, matching style of existing header cells (E1, etc.)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Timestamps for each data row (as text, to match inlineStr type in source data)
$timestamps = @(
    "2021-10-05 13:39:19.452871",
    "2021-10-05 13:39:19.452885",
    "2021-10-05 13:39:19.452889",
    "2021-10-05 13:39:19.452893",
    "2021-10-05 13:39:19.452896",
    "2021-10-05 13:39:19.452899",
    "2021-10-05 13:39:19.452902",
    "2021-10-05 13:39:19.452906",
    "2021-10-05 13:39:19.452909",
    "2021-10-05 13:39:19.452912",
    "2021-10-05 13:39:19.452915",
    "2021-10-05 13:39:19.452918"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
